$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell "D2" "26.863.07"
Set-TextCell "E2" "  +0.45%  "

Set-TextCell "D3" "1.814.43"
Set-TextCell "E3" "  +0.90%  "

Set-TextCell "E4" "  +0.06%  "

Set-TextCell "D5" "308.66"
Set-TextCell "E5" "  -0.01%  "

Set-TextCell "E6" "  +0.05%  "

Set-TextCell "D7" "0.4656"
Set-TextCell "E7" "  +1.01%  "

Set-TextCell "D8" "0.3679"
Set-TextCell "E8" "  -0.97%  "

Set-TextCell "E9" "  +1.44%  "

Set-TextCell "D10" "0.8692"

Set-TextCell "D11" "20.38"
Set-TextCell "E11" "  -0.04%  "

Set-TextCell "D12" "1.785.81"
Set-TextCell "E12" "  +1.00%  "

Set-TextCell "D13" "5.370"
Set-TextCell "E13" "  +1.14%  "

Set-TextCell "D14" "6.511"
Set-TextCell "E14" "  +0.34%  "

Set-TextCell "D15" "0.07054"
Set-TextCell "E15" "  +0.21%  "

Set-TextCell "D16" "91.58"
Set-TextCell "E16" "  +1.04%  "

Set-TextCell "E17" "  +0.13%  "

Set-TextCell "D18" "0.000008685"
Set-TextCell "E18" "  +0.55%  "

Set-TextCell "E19" "  +0.03%  "

Set-TextCell "D20" "14.73"
Set-TextCell "E20" "  +0.78%  "

Set-TextCell "D21" "26.900.11"
Set-TextCell "E21" "  +0.56%  "

Set-TextCell "D22" "5.334"
Set-TextCell "E22" "  +0.89%  "

Set-TextCell "D23" "10.56"
Set-TextCell "E23" "  -0.77%  "

Set-TextCell "D24" "2.053.84"
Set-TextCell "E24" "  +2.82%  "

Set-TextCell "D25" "1.896"
Set-TextCell "E25" "  -0.62%  "

Set-TextCell "D26" "150.32"
Set-TextCell "E26" "  +0.15%  "

Set-TextCell "D27" "2.172"
Set-TextCell "E27" "  +1.61%  "

Set-TextCell "E28" "  +1.01%  "

Set-TextCell "D29" "5.316"
Set-TextCell "E29" "  +1.91%  "

Set-TextCell "D30" "115.67"
Set-TextCell "E30" "  +1.33%  "

Set-TextCell "D31" "0.08921"
Set-TextCell "E31" "  +0.32%  "

Set-TextCell "D32" "0.7651"
Set-TextCell "E32" "  +1.38%  "

Set-TextCell "D33" "1.162"
Set-TextCell "E33" "  +0.35%  "

Set-TextCell "E34" "  +1.81%  "

Set-TextCell "E35" "  +0.52%  "

Set-TextCell "E37" "  -2.75%  "

Set-TextCell "B38" "Hedera"
Set-TextCell "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D38" "0.05283"
Set-TextCell "E38" "  +1.51%  "

Set-TextCell "B39" "VeChain"
Set-TextCell "C39" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D39" "0.01958"
Set-TextCell "E39" "  +0.91%  "

Set-TextCell "D40" "2.933"
Set-TextCell "E40" "  +1.22%  "

Set-TextCell "D41" "7.251"
Set-TextCell "E41" "  +1.26%  "

Set-TextCell "D42" "0.5321"
Set-TextCell "E42" "  +1.99%  "

Set-TextCell "D43" "2.332"
Set-TextCell "E43" "  -1.65%  "

Set-TextCell "D44" "0.1658"
Set-TextCell "E44" "  +0.73%  "

Set-TextCell "D45" "8.415"
Set-TextCell "E45" "  -0.93%  "

Set-TextCell "D46" "0.4916"
Set-TextCell "E46" "  -1.76%  "

Set-TextCell "D47" "10.39"
Set-TextCell "E47" "  +1.15%  "

Set-TextCell "D48" "1.001"
Set-TextCell "E48" "  +0.09%  "

Set-TextCell "D49" "1.669"
Set-TextCell "E49" "  +1.48%  "

Set-TextCell "D50" "103.69"
Set-TextCell "E50" "  -0.30%  "

Set-TextCell "D51" "0.06287"
Set-TextCell "E51" "  +0.10%  "
